# chore: update Sheets via scheduled runner
# Refreshes market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leve rows across all crafting-job sheets.
# Values are plain numbers (no formulas in this workbook), so cells are
# simply overwritten with the latest figures; a few cells that no longer
# have a computed profit value are cleared instead.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 2731.3333
$ws.Range("I19").Value = 949.5
$ws.Range("J19").Value = 3622.25
$ws.Range("K19").Value = 949.5
$ws.Range("L19").Value = 3622.25
$ws.Range("M19").Value = -774.5
$ws.Range("N19").Value = -3972.25
# row 32
$ws.Range("H32").Value = 1287.25
$ws.Range("I32").Value = 900
$ws.Range("K32").Value = 900
$ws.Range("M32").Value = -574
# row 43
$ws.Range("H43").Value = 1243.3684
$ws.Range("I43").Value = 399.5
$ws.Range("K43").Value = 399.5
$ws.Range("M43").Value = -330.5
# row 132
$ws.Range("H132").Value = 879.9761999999999
$ws.Range("I132").Value = 931.45715
$ws.Range("K132").Value = 2794.37145
$ws.Range("M132").Value = -264.3714499999996
# row 137
$ws.Range("H137").Value = 1841.6471
$ws.Range("I137").Value = 1257.2222
$ws.Range("K137").Value = 3771.6666
$ws.Range("M137").Value = -1221.6666
# row 138
$ws.Range("H138").Value = 2785.1562
$ws.Range("I138").Value = 3155.9
$ws.Range("J138").Value = 2167.25
$ws.Range("K138").Value = 9467.700000000001
$ws.Range("L138").Value = 6501.75
$ws.Range("M138").Value = -4327.700000000001
$ws.Range("N138").Value = -16781.75
# row 139
$ws.Range("H139").Value = 46949.5
$ws.Range("J139").Value = 46949.5
$ws.Range("L139").Value = 46949.5
$ws.Range("N139").Value = -57229.5
# row 141
$ws.Range("H141").Value = 2963.889
$ws.Range("I141").Value = 1496.375
$ws.Range("K141").Value = 4489.125
$ws.Range("M141").Value = 690.875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 5075.641
$ws.Range("I32").Value = 3135.6072
$ws.Range("J32").Value = 10013.909
$ws.Range("K32").Value = 3135.6072
$ws.Range("L32").Value = 10013.909
$ws.Range("M32").Value = -2848.6072
$ws.Range("N32").Value = -10587.909
# row 45
$ws.Range("H45").Value = 1655.1428
$ws.Range("J45").Value = 1789.2727
$ws.Range("L45").Value = 1789.2727
$ws.Range("N45").Value = -2543.2727
# row 88
$ws.Range("H88").Value = 3050.5
$ws.Range("I88").Value = 2350.875
$ws.Range("J88").Value = 4449.75
$ws.Range("K88").Value = 2350.875
$ws.Range("L88").Value = 4449.75
$ws.Range("M88").Value = -1944.875
$ws.Range("N88").Value = -5261.75
# row 91
$ws.Range("H91").Value = 3050.5
$ws.Range("I91").Value = 2350.875
$ws.Range("J91").Value = 4449.75
$ws.Range("K91").Value = 2350.875
$ws.Range("L91").Value = 4449.75
$ws.Range("M91").Value = -946.875
$ws.Range("N91").Value = -7257.75
# row 96
$ws.Range("H96").Value = 10172
$ws.Range("J96").Value = 10172
$ws.Range("L96").Value = 10172
$ws.Range("N96").Value = -15664
# row 102
$ws.Range("H102").Value = 2024.1818
$ws.Range("I102").Value = 1520.875
$ws.Range("K102").Value = 1520.875
$ws.Range("M102").Value = 101.125
# row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 119246.47
$ws.Range("I86").Value = 1449.375
$ws.Range("J86").Value = 223955
$ws.Range("K86").Value = 1449.375
$ws.Range("L86").Value = 223955
$ws.Range("M86").Value = -326.375
$ws.Range("N86").Value = -226201
# row 89
$ws.Range("H89").Value = 119246.47
$ws.Range("I89").Value = 1449.375
$ws.Range("J89").Value = 223955
$ws.Range("K89").Value = 7246.875
$ws.Range("L89").Value = 1119775
$ws.Range("M89").Value = -1630.875
$ws.Range("N89").Value = -1131007
# row 99
$ws.Range("H99").Value = 1330.75
$ws.Range("I99").Value = 1024.3334
$ws.Range("K99").Value = 1024.3334
$ws.Range("M99").Value = 473.6666
# row 134
$ws.Range("H134").Value = 12685.407
$ws.Range("I134").Value = 13092.174
$ws.Range("J134").Value = 10346.5
$ws.Range("K134").Value = 39276.522
$ws.Range("L134").Value = 31039.5
$ws.Range("M134").Value = -36741.522
$ws.Range("N134").Value = -36109.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2996.465
$ws.Range("I31").Value = 1399.619
$ws.Range("K31").Value = 1399.619
$ws.Range("M31").Value = -1104.619
# row 34
$ws.Range("H34").Value = 2996.465
$ws.Range("I34").Value = 1399.619
$ws.Range("K34").Value = 1399.619
$ws.Range("M34").Value = -1197.619
# row 93
$ws.Range("H93").Value = 12450
$ws.Range("I93").Value = 4900
$ws.Range("J93").Value = 20000
$ws.Range("K93").Value = 4900
$ws.Range("L93").Value = 20000
$ws.Range("M93").Value = -3028
$ws.Range("N93").Value = -23744
# row 99
$ws.Range("H99").Value = 2355.5557
$ws.Range("J99").Value = 2800
$ws.Range("L99").Value = 2800
$ws.Range("N99").Value = -5796
# row 126
$ws.Range("H126").Value = 2355.5557
$ws.Range("J126").Value = 2800
$ws.Range("L126").Value = 8400
$ws.Range("N126").Value = -13340

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1087
$ws.Range("N2").ClearContents()
# row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# row 107
$ws.Range("H107").Value = 786.625
$ws.Range("J107").Value = 786.625
$ws.Range("L107").Value = 2359.875
$ws.Range("N107").Value = -6199.875
# row 131
$ws.Range("H131").Value = 791.99
$ws.Range("I131").Value = 517.8333
$ws.Range("J131").Value = 809.4894
$ws.Range("K131").Value = 1553.4999
$ws.Range("L131").Value = 2428.4682
$ws.Range("M131").Value = 3486.5001
$ws.Range("N131").Value = -12508.4682

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1736.8334
$ws.Range("I97").Value = 1862
$ws.Range("J97").Value = 1111
$ws.Range("K97").Value = 1862
$ws.Range("L97").Value = 1111
$ws.Range("M97").Value = -1366
$ws.Range("N97").Value = -2103
# row 102
$ws.Range("H102").Value = 2803.2104
$ws.Range("I102").Value = 3431
$ws.Range("J102").Value = 2346.6365
$ws.Range("K102").Value = 3431
$ws.Range("L102").Value = 2346.6365
$ws.Range("M102").Value = -1809
$ws.Range("N102").Value = -5590.636500000001
# row 126
$ws.Range("H126").Value = 3144538.8
$ws.Range("I126").Value = 5053115
$ws.Range("K126").Value = 15159345
$ws.Range("M126").Value = -15156875
# row 132
$ws.Range("H132").Value = 2267055.8
$ws.Range("I132").Value = 3848710.5
$ws.Range("K132").Value = 11546131.5
$ws.Range("M132").Value = -11543601.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 3283.5
$ws.Range("I16").Value = 3427.4546
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 3427.4546
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = -3257.4546
$ws.Range("N16").Value = -2040
# row 40
$ws.Range("H40").Value = 9625.75
$ws.Range("I40").Value = 1502
$ws.Range("J40").Value = 12333.667
$ws.Range("K40").Value = 1502
$ws.Range("L40").Value = 12333.667
$ws.Range("M40").Value = -1366
$ws.Range("N40").Value = -12605.667
# row 55
$ws.Range("H55").Value = 374.18182
$ws.Range("I55").Value = 272.66666
$ws.Range("J55").Value = 496
$ws.Range("K55").Value = 272.66666
$ws.Range("L55").Value = 496
$ws.Range("M55").Value = -99.66665999999998
$ws.Range("N55").Value = -842
# row 132
$ws.Range("H132").Value = 2559.7
$ws.Range("I132").Value = 3499
$ws.Range("K132").Value = 10497
$ws.Range("M132").Value = -7967

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 1022
$ws.Range("I81").Value = 488
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 976
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = 85
$ws.Range("N81").Value = -4522
# row 84
$ws.Range("H84").Value = 1022
$ws.Range("I84").Value = 488
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 4880
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = 424
$ws.Range("N84").Value = -22608
# row 122
$ws.Range("H122").Value = 116359.25
$ws.Range("I122").Value = 154312.33
$ws.Range("K122").Value = 462936.99
$ws.Range("M122").Value = -460486.99
# row 132
$ws.Range("H132").Value = 9190.380999999999
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# row 136
$ws.Range("H136").Value = 61733090
$ws.Range("I136").Value = 111116350
$ws.Range("K136").Value = 333349050
$ws.Range("M136").Value = -333346500
